$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.5
$ws.Range("G3").Value = 2.45
$ws.Range("N6").Value = 1.9
$ws.Range("O6").Value = 1.95
$ws.Range("P6").Value = 1.4
$ws.Range("AA7").Value = 5.9
$ws.Range("AB7").Value = 15.5
$ws.Range("AC7").Value = 80
$ws.Range("AE7").Value = 7.3
$ws.Range("AF7").Value = 12.5
$ws.Range("AG7").Value = 10.25
$ws.Range("AH7").Value = 32
$ws.Range("AI7").Value = 26
$ws.Range("AJ7").Value = 40
$ws.Range("G7").Value = 2.57
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 2.7
$ws.Range("J7").Value = 1.09
$ws.Range("K7").Value = 6
$ws.Range("P7").Value = 1.5
$ws.Range("Q7").Value = 2.27
$ws.Range("R7").Value = 1.83
$ws.Range("S7").Value = 1.78
$ws.Range("U7").Value = 12
$ws.Range("V7").Value = 9.75
$ws.Range("W7").Value = 29
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 37
$ws.Range("Z7").Value = 7.6
$ws.Range("AA8").Value = 6
$ws.Range("AB8").Value = 19.5
$ws.Range("AC8").Value = 120
$ws.Range("AE8").Value = 7.5
$ws.Range("AF8").Value = 16.5
$ws.Range("AH8").Value = 55
$ws.Range("AI8").Value = 45
$ws.Range("G8").Value = 2.12
$ws.Range("H8").Value = 2.95
$ws.Range("I8").Value = 3.55
$ws.Range("L8").Value = 1.5
$ws.Range("M8").Value = 2.27
$ws.Range("N8").Value = 2.42
$ws.Range("O8").Value = 1.44
$ws.Range("R8").Value = 2.07
$ws.Range("S8").Value = 1.6
$ws.Range("T8").Value = 5.6
$ws.Range("U8").Value = 8.75
$ws.Range("V8").Value = 9.5
$ws.Range("W8").Value = 20
$ws.Range("X8").Value = 21
$ws.Range("Y8").Value = 45
$ws.Range("Z8").Value = 6.3
$ws.Range("G10").Value = 2.15
$ws.Range("J10").Value = 1.07
$ws.Range("K10").Value = 9
$ws.Range("O10").Value = 1.57
$ws.Range("P10").Value = 1.5
$ws.Range("S10").Value = 1.72
$ws.Range("AA11").Value = 5.7
$ws.Range("AB11").Value = 13
$ws.Range("AE11").Value = 9
$ws.Range("AF11").Value = 18.5
$ws.Range("AG11").Value = 11.75
$ws.Range("AH11").Value = 50
$ws.Range("AI11").Value = 32
$ws.Range("AJ11").Value = 37
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 4.2
$ws.Range("L11").Value = 1.33
$ws.Range("M11").Value = 3.05
$ws.Range("N11").Value = 1.98
$ws.Range("O11").Value = 1.65
$ws.Range("P11").Value = 1.38
$ws.Range("Q11").Value = 2.47
$ws.Range("U11").Value = 6.7
$ws.Range("V11").Value = 7
$ws.Range("W11").Value = 11.5
$ws.Range("X11").Value = 12
$ws.Range("Y11").Value = 23
$ws.Range("Z11").Value = 8.5
$ws.Range("AE12").Value = 6.3
$ws.Range("AF12").Value = 9.5
$ws.Range("AG12").Value = 8.25
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 17.5
$ws.Range("AJ12").Value = 27
$ws.Range("G12").Value = 2.62
$ws.Range("I12").Value = 2.45
$ws.Range("N12").Value = 2.07
$ws.Range("T12").Value = 6.5
$ws.Range("U12").Value = 10.25
$ws.Range("V12").Value = 8.5
$ws.Range("W12").Value = 23
$ws.Range("X12").Value = 19
$ws.Range("Y12").Value = 28
$ws.Range("AA13").Value = 6
$ws.Range("AB13").Value = 15
$ws.Range("AC13").Value = 75
$ws.Range("AE13").Value = 10.25
$ws.Range("AF13").Value = 24
$ws.Range("AG13").Value = 14.5
$ws.Range("AH13").Value = 80
$ws.Range("AI13").Value = 45
$ws.Range("AJ13").Value = 50
$ws.Range("G13").Value = 1.6
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 5.3
$ws.Range("L13").Value = 1.35
$ws.Range("M13").Value = 2.95
$ws.Range("N13").Value = 2.02
$ws.Range("O13").Value = 1.62
$ws.Range("P13").Value = 1.4
$ws.Range("Q13").Value = 2.42
$ws.Range("T13").Value = 4.85
$ws.Range("U13").Value = 5.7
$ws.Range("V13").Value = 7
$ws.Range("W13").Value = 9.25
$ws.Range("X13").Value = 11.75
$ws.Range("Y13").Value = 25
$ws.Range("Z13").Value = 8
$ws.Range("J14").Value = 1.05
$ws.Range("L14").Value = 1.33
$ws.Range("O14").Value = 1.63
$ws.Range("J15").Value = 1.07
$ws.Range("L15").Value = 1.47
$ws.Range("O15").Value = 1.47
$ws.Range("J16").Value = 1.08
$ws.Range("L16").Value = 1.5
$ws.Range("M16").Value = 2.37
$ws.Range("O16").Value = 1.41
$ws.Range("J17").Value = 1.05
$ws.Range("L17").Value = 1.3
$ws.Range("O17").Value = 1.67
$ws.Range("J18").Value = 1.07
$ws.Range("L18").Value = 1.41
$ws.Range("M18").Value = 2.62
$ws.Range("O18").Value = 1.5
$ws.Range("O19").Value = 1.67
$ws.Range("N20").Value = 2.07
$ws.Range("O20").Value = 1.69
$ws.Range("J21").Value = 1.1
$ws.Range("K21").Value = 7
$ws.Range("O21").Value = 1.47
$ws.Range("O22").Value = 1.44
$ws.Range("AE23").Value = 11
$ws.Range("L23").Value = 1.18
$ws.Range("M23").Value = 4.5
$ws.Range("N23").Value = 1.6
$ws.Range("O23").Value = 2.3
$ws.Range("P23").Value = 1.29
$ws.Range("Q23").Value = 3.5
$ws.Range("R23").Value = 1.53
$ws.Range("S23").Value = 2.38
$ws.Range("Y23").Value = 23
$ws.Range("K31").Value = 10
$ws.Range("P33").Value = 1.33
$ws.Range("AB36").Value = 16.5
$ws.Range("AC36").Value = 90
$ws.Range("AD36").Value = 800
$ws.Range("AE36").Value = 9.5
$ws.Range("AF36").Value = 19
$ws.Range("AG36").Value = 13
$ws.Range("AJ36").Value = 50
$ws.Range("H36").Value = 3.25
$ws.Range("I36").Value = 3.7
$ws.Range("L36").Value = 1.36
$ws.Range("M36").Value = 2.65
$ws.Range("N36").Value = 2.05
$ws.Range("O36").Value = 1.6
$ws.Range("P36").Value = 1.42
$ws.Range("Q36").Value = 2.47
$ws.Range("R36").Value = 1.88
$ws.Range("S36").Value = 1.72
$ws.Range("T36").Value = 6.2
$ws.Range("U36").Value = 8.5
$ws.Range("V36").Value = 8.75
$ws.Range("X36").Value = 17.5
$ws.Range("AA38").Value = 12
$ws.Range("AB38").Value = 15
$ws.Range("AC38").Value = 34
$ws.Range("AD38").Value = 81
$ws.Range("AE38").Value = 29
$ws.Range("AF38").Value = 41
$ws.Range("AG38").Value = 19
$ws.Range("AH38").Value = 51
$ws.Range("AI38").Value = 34
$ws.Range("AJ38").Value = 29
$ws.Range("J38").Value = 1.01
$ws.Range("K38").Value = 34
$ws.Range("L38").Value = 1.07
$ws.Range("M38").Value = 9
$ws.Range("P38").Value = 1.17
$ws.Range("Q38").Value = 5
$ws.Range("R38").Value = 1.4
$ws.Range("S38").Value = 2.75
$ws.Range("T38").Value = 15
$ws.Range("U38").Value = 12
$ws.Range("X38").Value = 11
$ws.Range("Y38").Value = 15
$ws.Range("Z38").Value = 34
